# regen sval data to filter save games
# Update the B:G numeric columns for rows 2-10 with the regenerated values
# (dates in column A and the Win flag in column F are unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.6545652718822623, 1.626987699542094, 3.223369029078222, 13.86384647080068, 19.36876847130326)
    3  = @(0.6545652718822623, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 6.038307959104277)
    4  = @(0.6545652718822623, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 6.038307959104277)
    5  = @(0.6545652718822623, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 6.038307959104277)
    6  = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    7  = @(3.272327238179451, 9.983522426115931, 18.71679738969934, 13.86384647080068, 45.8364935247954)
    8  = @(0.1169995834814548, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 2.426980108624251)
    9  = @(0.1169995834814548, 0.04103571897497393, 0.7210945179870265, 0.5333859586016987, 1.412515779045154)
    10 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B - TB
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C - d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D - K
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E - IP
    $ws.Cells.Item($row, 7).Value = $vals[4]  # G - sum
}
